$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "http://www.mendeley.com/c/5069088332/g/2058663/mlbsmammalspdf-applicationpdf-object/"
$ws.Range("B26").Value = "http://www.mendeley.com/c/4981987782/g/2058663/usfs-bartlett-experimental-forest/"
